{"js": "// Commit: \"atualizado o modelo do projeto + escopo para ser feito\"\n//\n// Net text-level edits inside the \"Amigo do Professor\" project idea:\n//   1) \"dificil\" -> \"Dif\u00edcil\"              (fix capitalization + accent)\n//   2) \"certo  e\" -> \"certo e\"             (collapse a stray double space)\n\nconst body = context.document.body;\n\n// 1) \"dificil\" -> \"Dif\u00edcil\"\nconst hard = body.search(\"dificil\", { matchCase: false });\nhard.load(\"items/text\");\nawait context.sync();\n\nif (hard.items.length > 0) {\n  hard.items[0].insertText(\"Dif\u00edcil\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"certo  e\" (two spaces) -> \"certo e\" (one space)\nconst doubleSpace = body.search(\"certo  e\", { matchCase: true });\ndoubleSpace.load(\"items/text\");\nawait context.sync();\n\nif (doubleSpace.items.length > 0) {\n  doubleSpace.items[0].insertText(\"certo e\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Commit: \"atualizado o modelo do projeto + escopo para ser feito\"\n#\n# Net text-level edits inside the \"Amigo do Professor\" project idea:\n#   1) \"dificil\" -> \"Dif\u00edcil\"   (fix capitalization + accent)\n#   2) \"certo  e\" -> \"certo e\"  (collapse a stray double space)\n\n$d = $word.ActiveDocument\n\n# 1) \"dificil\" -> \"Dif\u00edcil\"\n$find1 = $d.Content.Find\n$find1.Text = \"dificil\"\n$find1.Replacement.Text = \"Dif\u00edcil\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# 2) \"certo  e\" (two spaces) -> \"certo e\" (one space)\n$find2 = $d.Content.Find\n$find2.Text = \"certo  e\"\n$find2.Replacement.Text = \"certo e\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
